# Applies the change: replace "Mifos style" in ProductLoanInput!B17 with
# "Penalties, Fees, Interest, Principal order", update its style to a
# left/top aligned (no wrap) variant, and make B17 the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProductLoanInput")

$cell = $ws.Range("B17")
$cell.Value = "Penalties, Fees, Interest, Principal order"

$cell.HorizontalAlignment = -4131   # xlLeft
$cell.VerticalAlignment = -4160     # xlTop
$cell.WrapText = $false

$cell.Select()
